# Updated cryptos list on Sun Jul 30 15:34:28 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.457.17"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").Value = "1.883.64"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'0.7194"
$ws.Range("E5").Value = "  +1.58%  "

$ws.Range("D6").Value = "'243.67"
$ws.Range("E6").Value = "  +0.77%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "'0.07968"
$ws.Range("E8").Value = "  +2.11%  "

$ws.Range("D9").Value = "'0.3155"
$ws.Range("E9").Value = "  +1.53%  "

$ws.Range("D10").Value = "'25.07"
$ws.Range("E10").Value = "  +0.08%  "

$ws.Range("D11").Value = "'0.08141"
$ws.Range("E11").Value = "  -3.07%  "

$ws.Range("D12").Value = "1.900.53"
$ws.Range("E12").Value = "  +1.56%  "

$ws.Range("D13").Value = "'5.256"
$ws.Range("E13").Value = "  +0.29%  "

$ws.Range("E14").Value = "  +4.28%  "

$ws.Range("D15").Value = "'0.7121"
$ws.Range("E15").Value = "  -0.75%  "

$ws.Range("D16").Value = "'6.405"
$ws.Range("E16").Value = "  +4.43%  "

$ws.Range("D17").Value = "'0.000008449"
$ws.Range("E17").Value = "  +0.99%  "

$ws.Range("D18").Value = "29.463.21"
$ws.Range("E18").Value = "  +0.55%  "

$ws.Range("D19").Value = "'254.69"
$ws.Range("E19").Value = "  +5.94%  "

$ws.Range("E20").Value = "  +1.20%  "

$ws.Range("D21").Value = "2.139.76"
$ws.Range("E21").Value = "  +0.81%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").Value = "'7.800"
$ws.Range("E23").Value = "  +0.49%  "

$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("D25").Value = "'0.1591"
$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("D26").Value = "'9.092"
$ws.Range("E26").Value = "  +0.70%  "

$ws.Range("D27").Value = "'162.79"
$ws.Range("E27").Value = "  +0.07%  "

$ws.Range("D28").Value = "'19.09"
$ws.Range("E28").Value = "  +3.17%  "

$ws.Range("D29").Value = "'1.510"
$ws.Range("E29").Value = "  +0.45%  "

$ws.Range("D30").Value = "'4.431"
$ws.Range("E30").Value = "  +0.51%  "

$ws.Range("E31").Value = "  -0.62%  "

$ws.Range("D32").Value = "'1.224"
$ws.Range("E32").Value = "  -1.82%  "

$ws.Range("D33").Value = "'0.05327"
$ws.Range("E33").Value = "  -0.56%  "

$ws.Range("D34").Value = "'1.954"
$ws.Range("E34").Value = "  +0.84%  "

$ws.Range("D35").Value = "'0.7584"
$ws.Range("E35").Value = "  +1.23%  "

$ws.Range("E36").Value = "  +0.68%  "

$ws.Range("E37").Value = "  +0.71%  "

$ws.Range("D38").Value = "'0.01898"
$ws.Range("E38").Value = "  +1.17%  "

$ws.Range("D39").Value = "1.272.21"
$ws.Range("E39").Value = "  +2.46%  "

$ws.Range("D40").Value = "'2.769"
$ws.Range("E40").Value = "  +1.33%  "

$ws.Range("D41").Value = "'6.475"
$ws.Range("E41").Value = "  -0.63%  "

$ws.Range("D42").Value = "'113.04"
$ws.Range("E42").Value = "  +3.54%  "

$ws.Range("D43").Value = "'74.57"
$ws.Range("E43").Value = "  +3.05%  "

$ws.Range("D44").Value = "'0.9057"
$ws.Range("E44").Value = "  +1.59%  "

$ws.Range("E45").Value = "  +3.45%  "

$ws.Range("E46").Value = "  +0.12%  "

$ws.Range("D47").Value = "2.037.21"
$ws.Range("E47").Value = "  +0.65%  "

$ws.Range("E48").Value = "  +0.96%  "

$ws.Range("D49").Value = "'0.5203"
$ws.Range("E49").Value = "  +0.07%  "

$ws.Range("D50").Value = "'9.533"
$ws.Range("E50").Value = "  +1.07%  "

$ws.Range("D51").Value = "'0.4376"
$ws.Range("E51").Value = "  +0.87%  "
